$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.659.40"
$ws.Range("D2").Style = 'Normal'
$ws.Range("E2").Value = '  +5.02%  '

$ws.Range("D3").Value = "'3.624.89"
$ws.Range("D3").Style = 'Normal'
$ws.Range("E3").Value = '  +17.80%  '

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = 'Normal'
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").Value = "'593.23"
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  +3.06%  '

$ws.Range("D6").Value = "'186.12"
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  +8.66%  '

$ws.Range("D7").Value = "'3.623.20"
$ws.Range("D7").Style = 'Normal'
$ws.Range("E7").Value = '  +17.88%  '

$ws.Range("E8").Value = '  +0.05%  '

$ws.Range("D9").Value = "'0.534"
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = '  +4.80%  '

$ws.Range("E10").Value = '  +8.86%  '

$ws.Range("D11").Value = "'6.51"
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  +3.95%  '

$ws.Range("D12").Value = "'0.496"
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = '  +5.52%  '

$ws.Range("B13").Value = 'Avalanche'
$ws.Range("C13").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D13").Value = "'39.37"
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '  +9.84%  '

$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").Value = "'0.0000254"
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  +6.36%  '

$ws.Range("D15").Value = "'4.229.87"
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  +17.80%  '

$ws.Range("D16").Value = "'3.614.13"
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = '  +17.49%  '

$ws.Range("D17").Value = "'69.765.92"
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  +5.20%  '

$ws.Range("E18").Value = '  +2.04%  '

$ws.Range("D19").Value = "'7.52"
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  +8.12%  '

$ws.Range("D20").Value = "'17.29"
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  +4.26%  '

$ws.Range("D21").Value = "'509.31"
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  +4.59%  '

$ws.Range("D22").Value = "'9.14"
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  +19.32%  '

$ws.Range("E23").Value = '  +9.19%  '

$ws.Range("D24").Value = "'88.30"
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  +7.19%  '

$ws.Range("D25").Value = "'13.50"
$ws.Range("D25").Style = 'Normal'
$ws.Range("E25").Value = '  +7.07%  '

$ws.Range("D26").Value = "'2.40"
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  +8.61%  '

$ws.Range("D27").Value = "'10.77"
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  +6.78%  '

$ws.Range("E28").Value = '  +0.06%  '

$ws.Range("D29").Value = "'2.54"
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  +13.35%  '

$ws.Range("D30").Value = "'8.22"
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  +4.60%  '

$ws.Range("D31").Value = "'31.92"
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  +15.19%  '

$ws.Range("E32").Value = '  +6.08%  '

$ws.Range("E33").Value = '  +18.81%  '

$ws.Range("E34").Value = '  +5.45%  '

$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  -0.11%  '

$ws.Range("D36").Value = "'6.14"
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  +10.46%  '

$ws.Range("E37").Value = '  +7.26%  '

$ws.Range("D38").Value = "'0.335"
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  +11.21%  '

$ws.Range("D39").Value = "'46.96"
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  -1.89%  '

$ws.Range("E40").Value = '  +7.10%  '

$ws.Range("D41").Value = "'50.84"
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  +3.73%  '

$ws.Range("D42").Value = "'0.129"
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  +5.25%  '

$ws.Range("B43").Value = 'Cosmos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D43").Value = "'8.89"
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  +8.08%  '

$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = "'3.141.10"
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  +13.43%  '

$ws.Range("E45").Value = '  +9.57%  '

$ws.Range("D46").Value = "'399.91"
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  +9.73%  '

$ws.Range("D48").Value = "'27.89"
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  +15.06%  '

$ws.Range("D49").Value = "'136.06"
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  +0.98%  '

$ws.Range("E50").Value = '  +14.79%  '

$ws.Range("E51").Value = '  +0.02%  '
